# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1353530
$ws.Range("C4").Value = 6221
$ws.Range("D4").Value = 239157
$ws.Range("E4").Value = 1034022
$ws.Range("G4").Value = 314
$ws.Range("H4").Value = 80351

# Row 10 - Alemania
$ws.Range("B10").Value = 171704
$ws.Range("C10").Value = 380
$ws.Range("E10").Value = 19755
$ws.Range("F10").Value = 1581

# Row 11 - Brasil
$ws.Range("B11").Value = 156862
$ws.Range("C11").Value = 801
$ws.Range("E11").Value = 84438
$ws.Range("G11").Value = 83
$ws.Range("H11").Value = 10739

# Row 16 - India
$ws.Range("B16").Value = 67044
$ws.Range("C16").Value = 4236
$ws.Range("D16").Value = 20815
$ws.Range("E16").Value = 44022
$ws.Range("G16").Value = 106
$ws.Range("H16").Value = 2207

# Row 22 - Pakistan
$ws.Range("D22").Value = 8063
$ws.Range("E22").Value = 21612

# Row 24 - Ecuador
$ws.Range("B24").Value = 29559
$ws.Range("C24").Value = 488
$ws.Range("E24").Value = 23999
$ws.Range("G24").Value = 410
$ws.Range("H24").Value = 2127

# Row 25 - becomes Chile (new data), was Portugal
$ws.Range("A25").Value = "Chile"
$ws.Range("B25").Value = 28866
$ws.Range("C25").Value = 1647
$ws.Range("D25").Value = 13112
$ws.Range("E25").Value = 15442
$ws.Range("F25").Value = 544
$ws.Range("G25").Value = 8
$ws.Range("H25").Value = 312

# Row 26 - becomes Portugal (old Portugal data), was Chile
$ws.Range("A26").Value = "Portugal"
$ws.Range("B26").Value = 27581
$ws.Range("C26").Value = 175
$ws.Range("D26").Value = 2549
$ws.Range("E26").Value = 23897
$ws.Range("F26").Value = 112
$ws.Range("G26").Value = 9
$ws.Range("H26").Value = 1135

# Row 33 - Israel
$ws.Range("B33").Value = 16477
$ws.Range("C33").Value = 23
$ws.Range("D33").Value = 11430
$ws.Range("E33").Value = 4795
$ws.Range("G33").Value = 5
$ws.Range("H33").Value = 252
